# #fix: TROCA CV POR UMA
# Corrige quantidade (coluna F) e valor total (coluna H) de alguns itens
# da planilha de estoque, refletindo os numeros corretos.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Linha 24
$ws.Range("F24").Value = 1099.000
$ws.Range("H24").Value = 9995.35

# Linha 26
$ws.Range("F26").Value = 1285.000
$ws.Range("H26").Value = 13212.04

# Linha 39
$ws.Range("F39").Value = 1893.000
$ws.Range("H39").Value = 49308.35

# Linha 42
$ws.Range("F42").Value = 630.000
$ws.Range("H42").Value = 53532.46

# Linha 49
$ws.Range("F49").Value = 3283.000
$ws.Range("H49").Value = 12743.48

# Linha 76
$ws.Range("F76").Value = 25705.000
$ws.Range("H76").Value = 38852.12

# Linha 102
$ws.Range("F102").Value = 10620.300
$ws.Range("H102").Value = 24104.53

# Linha 111
$ws.Range("F111").Value = 2219.000
$ws.Range("H111").Value = 5570.48

# Linha 112 (somente a coluna H muda)
$ws.Range("H112").Value = 4800.50
